$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 776, shifting existing rows 776-817 down to 777-818.
$ws.Rows("776:776").Insert()

# Populate the newly inserted row with the new data point.
# Force column A to be stored as literal text (not auto-converted to a date
# serial number) by temporarily switching to a text number format, then
# clearing the format afterwards so the cell ends up with the same default
# (unstyled) appearance as every other date cell in the column.
$ws.Range("A776").NumberFormat = "@"
$ws.Range("A776").Value = "2026/02/09"
$ws.Range("A776").ClearFormats()

$ws.Range("B776").Value = "月"
$ws.Range("C776").Value = 19
$ws.Range("D776").Value = 109
